$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 639, pushing existing rows 639:680 down to 640:681
$ws.Rows.Item(639).Insert()

# Populate the newly inserted row 639 with the new weekly record
$ws.Cells.Item(639, 1).Value = 9
$ws.Cells.Item(639, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(639, 3).Value = "Metropolitana"
$ws.Cells.Item(639, 4).Value = 44706
$ws.Cells.Item(639, 5).Value = 13
$ws.Cells.Item(639, 6).Value = 100112023
$ws.Cells.Item(639, 7).Value = "Brócoli"
$ws.Cells.Item(639, 8).Value = "Sin especificar"
$ws.Cells.Item(639, 9).Value = "Primera"
$ws.Cells.Item(639, 10).Value = 4700
$ws.Cells.Item(639, 11).Value = 700
$ws.Cells.Item(639, 12).Value = 800
$ws.Cells.Item(639, 13).Value = 751
$ws.Cells.Item(639, 14).Value = "`$/unidad"
$ws.Cells.Item(639, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(639, 16).Value = 751
$ws.Cells.Item(639, 17).Value = 1
$ws.Cells.Item(639, 18).Value = "Hortaliza"
